# Set the relevant attendance cells from 0 to 1 on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cells = @("G3", "H3", "D4", "E4", "H5", "H6", "H7", "H8", "D9", "E9", "H10", "D11", "E11", "H12", "D13", "E13", "H14", "H15", "H16", "H17", "H18")

foreach ($cell in $cells) {
    $ws.Range($cell).Value = 1
}
